# Apply cryptos list update (prices + 1h volume %) generated from upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.570.36'
$ws.Cells.Item(2, 5).Value = '  -0.18%  '

$ws.Cells.Item(3, 4).Value = '1.853.83'
$ws.Cells.Item(3, 5).Value = '  -0.09%  '

$ws.Cells.Item(4, 4).Value = "'0.9976"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.71%  '

$ws.Cells.Item(5, 4).Value = "'265.19"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +2.52%  '

$ws.Cells.Item(6, 4).Value = "'0.9980"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.46%  '

$ws.Cells.Item(7, 4).Value = "'0.5225"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.27%  '

$ws.Cells.Item(8, 4).Value = "'0.3289"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.22%  '

$ws.Cells.Item(9, 4).Value = "'0.06827"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +1.12%  '

$ws.Cells.Item(10, 4).Value = "'18.90"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -2.89%  '

$ws.Cells.Item(11, 4).Value = "'0.7810"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.81%  '

$ws.Cells.Item(12, 4).Value = "'0.07782"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +1.03%  '

$ws.Cells.Item(13, 4).Value = '1.850.29'
$ws.Cells.Item(13, 5).Value = '  -0.76%  '

$ws.Cells.Item(14, 4).Value = "'88.56"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.61%  '

$ws.Cells.Item(15, 4).Value = "'5.029"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -0.65%  '

$ws.Cells.Item(16, 4).Value = "'0.9991"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.53%  '

$ws.Cells.Item(17, 4).Value = "'13.99"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.56%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = "'0.000007995"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +1.01%  '

$ws.Cells.Item(19, 2).Value = 'Dai'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(19, 4).Value = "'0.9994"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.42%  '

$ws.Cells.Item(20, 4).Value = '26.572.05'
$ws.Cells.Item(20, 5).Value = '  -0.40%  '

$ws.Cells.Item(21, 4).Value = '2.085.52'
$ws.Cells.Item(21, 5).Value = '  +0.60%  '

$ws.Cells.Item(22, 4).Value = "'4.653"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.76%  '

$ws.Cells.Item(23, 4).Value = "'9.556"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -2.01%  '

$ws.Cells.Item(24, 4).Value = "'6.004"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +0.07%  '

$ws.Cells.Item(25, 4).Value = "'144.67"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.21%  '

$ws.Cells.Item(26, 4).Value = "'2.240"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -5.26%  '

$ws.Cells.Item(27, 5).Value = '  -0.18%  '

$ws.Cells.Item(28, 4).Value = "'17.06"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.03%  '

$ws.Cells.Item(29, 4).Value = "'112.29"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +0.37%  '

$ws.Cells.Item(30, 4).Value = "'4.217"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.22%  '

$ws.Cells.Item(31, 4).Value = "'4.151"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -1.52%  '

$ws.Cells.Item(32, 4).Value = "'0.08768"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.13%  '

$ws.Cells.Item(33, 4).Value = "'0.04857"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.49%  '

$ws.Cells.Item(34, 4).Value = "'1.144"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.00%  '

$ws.Cells.Item(35, 4).Value = "'0.7207"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +1.36%  '

$ws.Cells.Item(36, 4).Value = "'2.853"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.29%  '

$ws.Cells.Item(37, 4).Value = "'3.104"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -1.68%  '

$ws.Cells.Item(38, 4).Value = "'0.01784"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -2.21%  '

$ws.Cells.Item(39, 4).Value = "'2.224"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -1.26%  '

$ws.Cells.Item(40, 4).Value = "'0.4911"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -1.31%  '

$ws.Cells.Item(41, 4).Value = "'0.9182"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +0.86%  '

$ws.Cells.Item(42, 4).Value = "'111.54"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -3.22%  '

$ws.Cells.Item(43, 4).Value = "'6.087"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.14%  '

$ws.Cells.Item(44, 4).Value = "'0.9981"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.38%  '

$ws.Cells.Item(45, 4).Value = "'7.775"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.12%  '

$ws.Cells.Item(46, 4).Value = "'0.4197"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -2.57%  '

$ws.Cells.Item(47, 4).Value = "'0.05953"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.46%  '

$ws.Cells.Item(48, 4).Value = "'9.086"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.47%  '

$ws.Cells.Item(49, 4).Value = "'0.1247"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -3.89%  '

$ws.Cells.Item(50, 4).Value = "'35.08"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.30%  '

$ws.Cells.Item(51, 4).Value = "'0.8940"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +2.97%  '
